$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.441378666666667
$ws.Range("H2").Value = 7.324135999999999
$ws.Range("I2").Value = 0.1119936059016048
$ws.Range("J2").Value = 0.1119936059016048
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.62804133333333
$ws.Range("N2").Value = 31.884124
$ws.Range("O2").Value = 0.1901422379349035
$ws.Range("P2").Value = 0.1901422379349035
$ws.Range("Q2").Value = 25.94707337965156
$ws.Range("R2").Value = 233.523660416864
$ws.Range("S2").Value = 0.02129471486053075
$ws.Range("T2").Value = 0.02129471486053075
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.441378666666667
$ws.Range("H3").Value = 7.324135999999999
$ws.Range("I3").Value = 0.1119936059016048
$ws.Range("J3").Value = 0.1119936059016048
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 26.47935433333333
$ws.Range("N3").Value = 79.438063
$ws.Range("O3").Value = 0.4737320390559845
$ws.Range("P3").Value = 0.4737320390559846
$ws.Range("Q3").Value = 64.64613077650755
$ws.Range("R3").Value = 581.8151769885679
$ws.Range("S3").Value = 0.05305495928499958
$ws.Range("T3").Value = 0.05305495928499959
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.441378666666667
$ws.Range("H4").Value = 7.324135999999999
$ws.Range("I4").Value = 0.1119936059016048
$ws.Range("J4").Value = 0.1119936059016048
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.626917
$ws.Range("N4").Value = 16.880751
$ws.Range("O4").Value = 0.1006690280454893
$ws.Range("P4").Value = 0.1006690280454894
$ws.Range("Q4").Value = 13.737435122904
$ws.Range("R4").Value = 123.636916106136
$ws.Range("S4").Value = 0.01127428745342413
$ws.Range("T4").Value = 0.01127428745342413
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.441378666666667
$ws.Range("H5").Value = 7.324135999999999
$ws.Range("I5").Value = 0.1119936059016048
$ws.Range("J5").Value = 0.1119936059016048
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.317432666666666
$ws.Range("N5").Value = 15.952298
$ws.Range("O5").Value = 0.09513216175939111
$ws.Range("P5").Value = 0.09513216175939114
$ws.Range("Q5").Value = 12.98186667383644
$ws.Range("R5").Value = 116.836800064528
$ws.Range("S5").Value = 0.01065419383264897
$ws.Range("T5").Value = 0.01065419383264897
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.441378666666667
$ws.Range("H6").Value = 7.324135999999999
$ws.Range("I6").Value = 0.1119936059016048
$ws.Range("J6").Value = 0.1119936059016048
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.84347
$ws.Range("N6").Value = 23.53041
$ws.Range("O6").Value = 0.1403245332042314
$ws.Range("P6").Value = 0.1403245332042315
$ws.Range("Q6").Value = 19.14888033064
$ws.Range("R6").Value = 172.33992297576
$ws.Range("S6").Value = 0.01571545047000135
$ws.Range("T6").Value = 0.01571545047000135
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.88577866666667
$ws.Range("H7").Value = 38.657336
$ws.Range("I7").Value = 0.5911106037886134
$ws.Range("J7").Value = 0.5911106037886134
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.62804133333333
$ws.Range("N7").Value = 31.884124
$ws.Range("O7").Value = 0.1901422379349035
$ws.Range("P7").Value = 0.1901422379349035
$ws.Range("Q7").Value = 136.9505882815182
$ws.Range("R7").Value = 1232.555294533664
$ws.Range("S7").Value = 0.112395093071419
$ws.Range("T7").Value = 0.112395093071419
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.88577866666667
$ws.Range("H8").Value = 38.657336
$ws.Range("I8").Value = 0.5911106037886134
$ws.Range("J8").Value = 0.5911106037886134
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 26.47935433333333
$ws.Range("N8").Value = 79.438063
$ws.Range("O8").Value = 0.4737320390559845
$ws.Range("P8").Value = 0.4737320390559846
$ws.Range("Q8").Value = 341.2070991755742
$ws.Range("R8").Value = 3070.863892580168
$ws.Range("S8").Value = 0.280028031640394
$ws.Range("T8").Value = 0.280028031640394
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.88577866666667
$ws.Range("H9").Value = 38.657336
$ws.Range("I9").Value = 0.5911106037886134
$ws.Range("J9").Value = 0.5911106037886134
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.626917
$ws.Range("N9").Value = 16.880751
$ws.Range("O9").Value = 0.1006690280454893
$ws.Range("P9").Value = 0.1006690280454894
$ws.Range("Q9").Value = 72.50720703770399
$ws.Range("R9").Value = 652.5648633393361
$ws.Range("S9").Value = 0.05950652995078206
$ws.Range("T9").Value = 0.05950652995078207
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.88577866666667
$ws.Range("H10").Value = 38.657336
$ws.Range("I10").Value = 0.5911106037886134
$ws.Range("J10").Value = 0.5911106037886134
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.317432666666666
$ws.Range("N10").Value = 15.952298
$ws.Range("O10").Value = 0.09513216175939111
$ws.Range("P10").Value = 0.09513216175939114
$ws.Range("Q10").Value = 68.51926041756977
$ws.Range("R10").Value = 616.673343758128
$ws.Range("S10").Value = 0.05623362957730971
$ws.Range("T10").Value = 0.05623362957730973
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.88577866666667
$ws.Range("H11").Value = 38.657336
$ws.Range("I11").Value = 0.5911106037886134
$ws.Range("J11").Value = 0.5911106037886134
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.84347
$ws.Range("N11").Value = 23.53041
$ws.Range("O11").Value = 0.1403245332042314
$ws.Range("P11").Value = 0.1403245332042315
$ws.Range("Q11").Value = 101.06921839864
$ws.Range("R11").Value = 909.6229655877601
$ws.Range("S11").Value = 0.08294731954870857
$ws.Range("T11").Value = 0.08294731954870858
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.472111000000001
$ws.Range("H12").Value = 19.416333
$ws.Range("I12").Value = 0.2968957903097819
$ws.Range("J12").Value = 0.2968957903097818
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.62804133333333
$ws.Range("N12").Value = 31.884124
$ws.Range("O12").Value = 0.1901422379349035
$ws.Range("P12").Value = 0.1901422379349035
$ws.Range("Q12").Value = 68.78586322192135
$ws.Range("R12").Value = 619.0727689972921
$ws.Range("S12").Value = 0.05645243000295375
$ws.Range("T12").Value = 0.05645243000295375
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.472111000000001
$ws.Range("H13").Value = 19.416333
$ws.Range("I13").Value = 0.2968957903097819
$ws.Range("J13").Value = 0.2968957903097818
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 26.47935433333333
$ws.Range("N13").Value = 79.438063
$ws.Range("O13").Value = 0.4737320390559845
$ws.Range("P13").Value = 0.4737320390559846
$ws.Range("Q13").Value = 171.3773204536644
$ws.Range("R13").Value = 1542.395884082979
$ws.Range("S13").Value = 0.140649048130591
$ws.Range("T13").Value = 0.140649048130591
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.472111000000001
$ws.Range("H14").Value = 19.416333
$ws.Range("I14").Value = 0.2968957903097819
$ws.Range("J14").Value = 0.2968957903097818
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.626917
$ws.Range("N14").Value = 16.880751
$ws.Range("O14").Value = 0.1006690280454893
$ws.Range("P14").Value = 0.1006690280454894
$ws.Range("Q14").Value = 36.41803141178701
$ws.Range("R14").Value = 327.7622827060831
$ws.Range("S14").Value = 0.02988821064128315
$ws.Range("T14").Value = 0.02988821064128315
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.472111000000001
$ws.Range("H15").Value = 19.416333
$ws.Range("I15").Value = 0.2968957903097819
$ws.Range("J15").Value = 0.2968957903097818
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.317432666666666
$ws.Range("N15").Value = 15.952298
$ws.Range("O15").Value = 0.09513216175939111
$ws.Range("P15").Value = 0.09513216175939114
$ws.Range("Q15").Value = 34.41501445369267
$ws.Range("R15").Value = 309.735130083234
$ws.Range("S15").Value = 0.02824433834943243
$ws.Range("T15").Value = 0.02824433834943244
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.472111000000001
$ws.Range("H16").Value = 19.416333
$ws.Range("I16").Value = 0.2968957903097819
$ws.Range("J16").Value = 0.2968957903097818
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.84347
$ws.Range("N16").Value = 23.53041
$ws.Range("O16").Value = 0.1403245332042314
$ws.Range("P16").Value = 0.1403245332042315
$ws.Range("Q16").Value = 50.76380846517001
$ws.Range("R16").Value = 456.87427618653
$ws.Range("S16").Value = 0.04166176318552152
$ws.Range("T16").Value = 0.04166176318552152
